$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 updates
# ---------------------------------------------------------------------------
# B2 holds a textual "2" (it was textual "1" before) - force as text via a
# formula that evaluates to a string, then convert the cell to a static value
# so it ends up as a plain shared-string cell (same representation Excel uses
# for the other text cells in this sheet).
$ws.Range("B2").Formula = "=""2"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("K2").Value = -10.53
$ws.Range("U2").Value = 4.459
$ws.Range("V2").Value = 0.1202210838500943
$ws.Range("W2").Value = -0.6977211815356541
$ws.Range("X2").Value = 0.1043418470071959
$ws.Range("Y2").Value = -0.8020630285428501
$ws.Range("AA2").Value = -0.5803107223600318
$ws.Range("AB2").Value = 0.10267079702392
$ws.Range("AC2").Value = -0.6829815193839519
$ws.Range("AD2").Value = 0.654
$ws.Range("AF2").Value = 0.654
$ws.Range("AG2").Value = -3.805
$ws.Range("AH2").Value = 0.01732725731242052
$ws.Range("AI2").Value = 0.01884726224783861
$ws.Range("AJ2").Value = -0.1143157578488809
$ws.Range("AK2").Value = -0.1258225587778182
$ws.Range("AL2").Value = 0.109
$ws.Range("AM2").Value = 0.109
$ws.Range("AN2").Value = 1.006153846153846
$ws.Range("AO2").Value = -93.85321100917432
$ws.Range("AP2").Value = -5.85384615384615
$ws.Range("AQ2").Value = -93.85321100917432

# ---------------------------------------------------------------------------
# Row 3 updates
# ---------------------------------------------------------------------------
$ws.Range("K3").Value = -5.18
$ws.Range("U3").Value = 3.94
$ws.Range("V3").Value = 0.123125
$ws.Range("W3").Value = -0.1598765432098765
$ws.Range("X3").Value = 0.1009231092645841
$ws.Range("Y3").Value = -0.2607996524744606
$ws.Range("AA3").Value = -0.1800247411940882
$ws.Range("AB3").Value = 0.1007722604069783
$ws.Range("AC3").Value = -0.2807970016010665
$ws.Range("AD3").Value = 0.126
$ws.Range("AF3").Value = 0.126
$ws.Range("AG3").Value = -3.814
$ws.Range("AH3").Value = 0.0039220569009525
$ws.Range("AI3").Value = 0.003607627555402852
$ws.Range("AJ3").Value = -0.135315404810899
$ws.Range("AK3").Value = -0.123087846124056
$ws.Range("AL3").Value = 0.019
$ws.Range("AM3").Value = 0.019
$ws.Range("AN3").Value = -0.05431034482758621
$ws.Range("AO3").Value = -291.0526315789474
$ws.Range("AP3").Value = 1.643965517241379
$ws.Range("AQ3").Value = -291.0526315789474

# ---------------------------------------------------------------------------
# New row 4
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Mongolia"
$ws.Range("B4").Value = "Aranjin Resources Ltd. (TSXV:ARJN)"
$ws.Range("C4").Value = "Metals & Mining"

$ws.Range("K4").Value = -5.35
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.519
$ws.Range("V4").Value = 0.1019646365422397
$ws.Range("W4").Value = -1.235565819861432
$ws.Range("X4").Value = 0.1077605847498077
$ws.Range("Y4").Value = -1.343326404611239
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = -0.9805967035259754
$ws.Range("AB4").Value = 0.1045693336408618
$ws.Range("AC4").Value = -1.085166037166837
$ws.Range("AD4").Value = 0.528
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.528
$ws.Range("AG4").Value = 0.009000000000000008
$ws.Range("AH4").Value = 0.09398362406550374
$ws.Range("AI4").Value = -2.336283185840708
$ws.Range("AJ4").Value = 0.001765051970974702
$ws.Range("AK4").Value = -0.01208053691275169
$ws.Range("AL4").Value = 0.09
$ws.Range("AM4").Value = 0.09
$ws.Range("AN4").Value = 0.1777777777777778
$ws.Range("AO4").Value = -52.22222222222223
$ws.Range("AP4").Value = 0.003030303030303033
$ws.Range("AQ4").Value = -52.22222222222223
